$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Coin name) updates ---
$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 2).Value = 'Algorand'
$ws.Cells.Item(51, 2).Value = 'WOONetwork'

# --- Column C (Link) updates ---
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'

# --- Column D (Price) updates ---
# Cells whose new value would be auto-parsed as a number need the column
# forced to Text format first so they are stored as literal strings,
# matching the source data (inline strings in the original workbook).
$ws.Cells.Item(2, 4).Value = '31.007.83'
$ws.Cells.Item(3, 4).Value = '1.954.31'
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '245.66'
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.001'
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4904'
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.2964'
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.06833'
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '19.19'
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '107.46'
$ws.Cells.Item(12, 4).Value = '1.955.09'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.07796'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.453'
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.7059'
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '282.90'
$ws.Cells.Item(17, 4).Value = '31.039.64'
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '13.23'
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.000007703'
$ws.Cells.Item(20, 4).Value = '2.202.93'
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.489'
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '1.001'
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '6.486'
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '9.854'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '169.85'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '20.02'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.206'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '0.1058'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.419'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.582'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.610'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '4.450'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.04952'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.7642'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.171'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.732'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.02012'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.702'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '6.615'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '2.139'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '74.15'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.4487'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '109.65'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.8843'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '8.154'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.001'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '988.04'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '9.380'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.1265'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.2576'

# --- Column E (Volume 1h) updates ---
$ws.Cells.Item(2, 5).Value = '  +1.01%  '
$ws.Cells.Item(3, 5).Value = '  -0.53%  '
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$ws.Cells.Item(5, 5).Value = '  -1.53%  '
$ws.Cells.Item(6, 5).Value = '  -0.07%  '
$ws.Cells.Item(7, 5).Value = '  +1.48%  '
$ws.Cells.Item(8, 5).Value = '  +0.35%  '
$ws.Cells.Item(9, 5).Value = '  +0.34%  '
$ws.Cells.Item(10, 5).Value = '  -1.28%  '
$ws.Cells.Item(11, 5).Value = '  -3.33%  '
$ws.Cells.Item(12, 5).Value = '  -0.12%  '
$ws.Cells.Item(13, 5).Value = '  +0.68%  '
$ws.Cells.Item(14, 5).Value = '  -0.79%  '
$ws.Cells.Item(15, 5).Value = '  +1.70%  '
$ws.Cells.Item(16, 5).Value = '  -3.95%  '
$ws.Cells.Item(17, 5).Value = '  +1.06%  '
$ws.Cells.Item(18, 5).Value = '  -1.03%  '
$ws.Cells.Item(19, 5).Value = '  -0.07%  '
$ws.Cells.Item(20, 5).Value = '  -0.42%  '
$ws.Cells.Item(21, 5).Value = '  -0.15%  '
$ws.Cells.Item(22, 5).Value = '  -3.20%  '
$ws.Cells.Item(23, 5).Value = '  +0.11%  '
$ws.Cells.Item(24, 5).Value = '  -2.63%  '
$ws.Cells.Item(25, 5).Value = '  -0.22%  '
$ws.Cells.Item(26, 5).Value = '  -0.02%  '
$ws.Cells.Item(27, 5).Value = '  -0.88%  '
$ws.Cells.Item(28, 5).Value = '  -0.11%  '
$ws.Cells.Item(29, 5).Value = '  -1.81%  '
$ws.Cells.Item(30, 5).Value = '  -1.50%  '
$ws.Cells.Item(31, 5).Value = '  -1.48%  '
$ws.Cells.Item(32, 5).Value = '  -1.24%  '
$ws.Cells.Item(33, 5).Value = '  -0.09%  '
$ws.Cells.Item(34, 5).Value = '  -3.04%  '
$ws.Cells.Item(35, 5).Value = '  -2.24%  '
$ws.Cells.Item(36, 5).Value = '  -0.88%  '
$ws.Cells.Item(37, 5).Value = '  -0.23%  '
$ws.Cells.Item(38, 5).Value = '  -2.66%  '
$ws.Cells.Item(39, 5).Value = '  -0.59%  '
$ws.Cells.Item(40, 5).Value = '  +8.42%  '
$ws.Cells.Item(41, 5).Value = '  +3.00%  '
$ws.Cells.Item(42, 5).Value = '  +5.58%  '
$ws.Cells.Item(43, 5).Value = '  +0.09%  '
$ws.Cells.Item(44, 5).Value = '  -1.86%  '
$ws.Cells.Item(45, 5).Value = '  +0.84%  '
$ws.Cells.Item(46, 5).Value = '  +9.35%  '
$ws.Cells.Item(47, 5).Value = '  -0.19%  '
$ws.Cells.Item(48, 5).Value = '  +8.32%  '
$ws.Cells.Item(49, 5).Value = '  +0.06%  '
$ws.Cells.Item(50, 5).Value = '  -1.48%  '
$ws.Cells.Item(51, 5).Value = '  +2.13%  '
